$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append "ТС-АВТО" test-case references to the two acceptance rows that
# cover the "publish" (auth_5) and "edit" (auth_6) scenarios.
$ws.Range("B13").Value = "ТС-ПУБ1, ТС-ПУБ2, ТС-ПУБ3, ТС-АВТО-1 "
$ws.Range("B14").Value = "ТС-РЕД1, ТС-РЕД2, ТС-РЕД3,  ТС-АВТО-12"

# Update view state: move selection near the edit.
$ws.Range("F14").Select()
